{"js": "const replacements = [\n  [\"2024-04-04 Thursday\", \"2024-04-05 Friday\"],\n  [\"719\u00d78=\", \"933\u00d74=\"],\n  [\"128\u00d75=\", \"738\u00d79=\"],\n  [\"631\u00d76=\", \"161\u00d74=\"],\n  [\"447\u00d79=\", \"518\u00d76=\"],\n  [\"802\u00d76=\", \"567\u00d79=\"],\n  [\"252\u00d76=\", \"841\u00d79=\"],\n  [\"482\u00d78=\", \"283\u00d77=\"],\n  [\"944\u00d77=\", \"359\u00d75=\"],\n  [\"900\u00d79=\", \"256\u00d74=\"],\n  [\"324\u00d76=\", \"171\u00d79=\"],\n  [\"401\u00d76=\", \"997\u00d75=\"],\n  [\"755\u00d73=\", \"600\u00d72=\"],\n  [\"169\u00d77=\", \"906\u00d78=\"],\n  [\"584\u00d78=\", \"281\u00d76=\"],\n  [\"750\u00d77=\", \"685\u00d79=\"],\n  [\"299\u00d73=\", \"378\u00d72=\"],\n  [\"837\u00d76=\", \"942\u00d74=\"],\n  [\"828\u00d75=\", \"955\u00d76=\"],\n  [\"373\u00d79=\", \"350\u00d79=\"],\n  [\"797\u00d78=\", \"913\u00d75=\"],\n  [\"463\u00d75=\", \"939\u00d74=\"],\n  [\"541\u00d74=\", \"224\u00d78=\"],\n  [\"209\u00d76=\", \"998\u00d75=\"],\n  [\"670\u00d78=\", \"366\u00d77=\"],\n  [\"904\u00d74=\", \"304\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-04 Thursday\", \"2024-04-05 Friday\"),\n    @(\"719\u00d78=\", \"933\u00d74=\"),\n    @(\"128\u00d75=\", \"738\u00d79=\"),\n    @(\"631\u00d76=\", \"161\u00d74=\"),\n    @(\"447\u00d79=\", \"518\u00d76=\"),\n    @(\"802\u00d76=\", \"567\u00d79=\"),\n    @(\"252\u00d76=\", \"841\u00d79=\"),\n    @(\"482\u00d78=\", \"283\u00d77=\"),\n    @(\"944\u00d77=\", \"359\u00d75=\"),\n    @(\"900\u00d79=\", \"256\u00d74=\"),\n    @(\"324\u00d76=\", \"171\u00d79=\"),\n    @(\"401\u00d76=\", \"997\u00d75=\"),\n    @(\"755\u00d73=\", \"600\u00d72=\"),\n    @(\"169\u00d77=\", \"906\u00d78=\"),\n    @(\"584\u00d78=\", \"281\u00d76=\"),\n    @(\"750\u00d77=\", \"685\u00d79=\"),\n    @(\"299\u00d73=\", \"378\u00d72=\"),\n    @(\"837\u00d76=\", \"942\u00d74=\"),\n    @(\"828\u00d75=\", \"955\u00d76=\"),\n    @(\"373\u00d79=\", \"350\u00d79=\"),\n    @(\"797\u00d78=\", \"913\u00d75=\"),\n    @(\"463\u00d75=\", \"939\u00d74=\"),\n    @(\"541\u00d74=\", \"224\u00d78=\"),\n    @(\"209\u00d76=\", \"998\u00d75=\"),\n    @(\"670\u00d78=\", \"366\u00d77=\"),\n    @(\"904\u00d74=\", \"304\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n\n$d.Save()\n"}
